$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the existing "foo" note (F2) to a more descriptive note about the
#    new null/empty-cell check. Since F2 is the sole cell referencing that
#    shared string, this updates the shared string table entry in place.
$ws.Range("F2").Value = "catching null values"

# 2) Add a new data row (row 3) reusing the same note text (so it reuses the
#    same shared string slot created above).
$ws.Range("F3").Value = "catching null values"

# 3) Simulate "empty" cells in row 2 by filling them with a single space, to
#    test the Excel-side check for null/empty values.
$ws.Range("B2:G2").Value = " "

# 4) Fill in the rest of the new row 3 with real data.
$ws.Range("B3").Value = 44652
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 4

# Copy the date format from B2 onto B3 so it keeps the same date style.
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)

# 5) Update the current selection to match where the user ended up.
$ws.Range("E11").Select()
